$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repulled data / mean calculation
$ws.Range("F2").Value = -9
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = 3
$ws.Range("F5").Value = -5
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = 1
